$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8
$ws.Range("B8").Value = "6_bert_uncased_vs_cased"
$ws.Range("C8").Value = 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = "29min"
$ws.Range("G8").Value = 0.9351
$ws.Range("H8").Value = 0.9
$ws.Range("I8").Value = "bert-uncased"
$ws.Range("J8").Value = "Trained on a 60/20/20 split"

# Row 9
$ws.Range("B9").Value = "6_bert_uncased_vs_cased"
$ws.Range("C9").Value = 8
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = "29min"
$ws.Range("G9").Value = 0.93
$ws.Range("I9").Value = "bert-cased"
$ws.Range("J9").Value = "Trained on a 60/20/20 split"
